# Update the workbook to rename the "get_started" / "out_of_scope" IDs
# used on the onboarding rows to the new "phrase_presentation" /
# "phrase_hors_sujet" IDs, and move the active selection to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "phrase_presentation"
$ws.Range("A3").Value = "phrase_presentation"
$ws.Range("A4").Value = "phrase_hors_sujet"

$ws.Range("C3").Select()
